# Autonomous_temporary: insert two new signal rows into the VCU_IGN_R2D
# message block (R2D_button_raw, Ignition_switch_raw), which pushes the
# ACU_status message block down by two rows, then append a brand new
# VCU_APPS_RAW message block at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autonomous_temporary")

# --- 1. Make room: shift the ACU_status block (rows 82-86) down to 84-88 ---
# (the plain data-row template "assi_state" - currently row 85 - shifts to
# row 87 once the insert below runs; row 1/2 message+header templates are
# above the insert point so they keep their row numbers)
$ws.Rows("82:83").Insert()

# --- 2. Populate the two new signal rows for VCU_IGN_R2D (uses the plain
#        data-row style, cloned from the existing "assi_state" row - now at
#        row 87 post-insert - so borders/fonts match exactly and no new
#        style entries are created) ---
$ws.Range("A87:K87").Copy()

$ws.Range("A81:K81").PasteSpecial(-4122)
$ws.Cells.Item(81, 1).Value = "R2D_button_raw"
$ws.Cells.Item(81, 2).Value = 48
$ws.Cells.Item(81, 3).Value = 8
$ws.Cells.Item(81, 4).Value = "Intel"
$ws.Cells.Item(81, 5).Value = $false
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = 0

$ws.Range("A82:K82").PasteSpecial(-4122)
$ws.Cells.Item(82, 1).Value = "Ignition_switch_raw"
$ws.Cells.Item(82, 2).Value = 56
$ws.Cells.Item(82, 3).Value = 8
$ws.Cells.Item(82, 4).Value = "Intel"
$ws.Cells.Item(82, 5).Value = $false
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 7).Value = 0

# row 83 stays blank - it is the separator before "Message: ACU_status"
# (now at row 84, shifted automatically by the Insert above).

# --- 3. Append the new VCU_APPS_RAW message block at the bottom ---
# row 89 stays blank - separator after the ACU_status block.

# Message header row (cloned from row 1's "Message: ACU_MS" style)
$ws.Range("A1:C1").Copy()
$ws.Range("A90:C90").PasteSpecial(-4122)
$ws.Cells.Item(90, 1).Value = "Message: VCU_APPS_RAW"
$ws.Cells.Item(90, 2).Value = "ID: 0x610"
$ws.Cells.Item(90, 3).Value = "Sender(s): VCU"

# Column header row (cloned from row 2's header style)
$ws.Range("A2:K2").Copy()
$ws.Range("A91:K91").PasteSpecial(-4122)
$ws.Cells.Item(91, 1).Value = "Signal Name"
$ws.Cells.Item(91, 2).Value = "Start Bit"
$ws.Cells.Item(91, 3).Value = "Length (bits)"
$ws.Cells.Item(91, 4).Value = "Byte Order"
$ws.Cells.Item(91, 5).Value = "Signed"
$ws.Cells.Item(91, 6).Value = "Factor"
$ws.Cells.Item(91, 7).Value = "Offset"
$ws.Cells.Item(91, 8).Value = "Min"
$ws.Cells.Item(91, 9).Value = "Max"
$ws.Cells.Item(91, 10).Value = "Unit"
$ws.Cells.Item(91, 11).Value = "Choices"

# Data rows (cloned from the plain data-row style again)
$ws.Range("A87:K87").Copy()

$ws.Range("A92:K92").PasteSpecial(-4122)
$ws.Cells.Item(92, 1).Value = "APPS_1_raw_bits"
$ws.Cells.Item(92, 2).Value = 0
$ws.Cells.Item(92, 3).Value = 8
$ws.Cells.Item(92, 4).Value = "Intel"
$ws.Cells.Item(92, 5).Value = $false
$ws.Cells.Item(92, 6).Value = 1
$ws.Cells.Item(92, 7).Value = 0

$ws.Range("A93:K93").PasteSpecial(-4122)
$ws.Cells.Item(93, 1).Value = "APPS_2_raw_bits"
$ws.Cells.Item(93, 2).Value = 8
$ws.Cells.Item(93, 3).Value = 8
$ws.Cells.Item(93, 4).Value = "Intel"
$ws.Cells.Item(93, 5).Value = $false
$ws.Cells.Item(93, 6).Value = 1
$ws.Cells.Item(93, 7).Value = 0

$ws.Range("A94:K94").PasteSpecial(-4122)
$ws.Cells.Item(94, 1).Value = "delta_raw"
$ws.Cells.Item(94, 2).Value = 16
$ws.Cells.Item(94, 3).Value = 8
$ws.Cells.Item(94, 4).Value = "Intel"
$ws.Cells.Item(94, 5).Value = $false
$ws.Cells.Item(94, 6).Value = 1
$ws.Cells.Item(94, 7).Value = 0

$ws.Range("A95:K95").PasteSpecial(-4122)
$ws.Cells.Item(95, 1).Value = "cpu_temp"
$ws.Cells.Item(95, 2).Value = 24
$ws.Cells.Item(95, 3).Value = 8
$ws.Cells.Item(95, 4).Value = "Intel"
$ws.Cells.Item(95, 5).Value = $false
$ws.Cells.Item(95, 6).Value = 1
$ws.Cells.Item(95, 7).Value = 0

Write-Output "edit applied"
